# Update the two-digit / one-digit division answer table with a new
# generated set of problems. The document contains a single table whose
# "problem" rows are 1, 5, 9, 13 and 17 (1-based), each with 5 columns.
# We address every cell explicitly via Table.Cell(row, col) so that the
# duplicate "73÷2=36, 1" values (row 1 col 3 and row 9 col 1, which must
# become two different strings) are each updated correctly.

$d = $word.ActiveDocument
$t = $d.Tables.Item(1)

$updates = @(
    @{ Row = 1;  Col = 1; New = "39÷3=13, 0" },
    @{ Row = 1;  Col = 2; New = "48÷5=9, 3" },
    @{ Row = 1;  Col = 3; New = "41÷5=8, 1" },
    @{ Row = 1;  Col = 4; New = "30÷6=5, 0" },
    @{ Row = 1;  Col = 5; New = "80÷6=13, 2" },

    @{ Row = 5;  Col = 1; New = "42÷4=10, 2" },
    @{ Row = 5;  Col = 2; New = "55÷6=9, 1" },
    @{ Row = 5;  Col = 3; New = "20÷7=2, 6" },
    @{ Row = 5;  Col = 4; New = "47÷2=23, 1" },
    @{ Row = 5;  Col = 5; New = "82÷8=10, 2" },

    @{ Row = 9;  Col = 1; New = "90÷9=10, 0" },
    @{ Row = 9;  Col = 2; New = "42÷2=21, 0" },
    @{ Row = 9;  Col = 3; New = "53÷8=6, 5" },
    @{ Row = 9;  Col = 4; New = "55÷4=13, 3" },
    @{ Row = 9;  Col = 5; New = "22÷7=3, 1" },

    @{ Row = 13; Col = 1; New = "83÷6=13, 5" },
    @{ Row = 13; Col = 2; New = "46÷7=6, 4" },
    @{ Row = 13; Col = 3; New = "75÷6=12, 3" },
    @{ Row = 13; Col = 4; New = "42÷3=14, 0" },
    @{ Row = 13; Col = 5; New = "24÷8=3, 0" },

    @{ Row = 17; Col = 1; New = "42÷5=8, 2" },
    @{ Row = 17; Col = 2; New = "44÷9=4, 8" },
    @{ Row = 17; Col = 3; New = "86÷2=43, 0" },
    @{ Row = 17; Col = 4; New = "33÷5=6, 3" },
    @{ Row = 17; Col = 5; New = "26÷7=3, 5" }
)

foreach ($u in $updates) {
    $cell = $t.Cell($u.Row, $u.Col)
    $cell.Range.Text = $u.New
}

Write-Host "Updated $($updates.Count) cells"
